# Auto-generated Excel COM-interop edit script
# Updates market-price / profit columns (H:N) on each crafting-job sheet
# to reflect a refreshed Universalis price pull (per the scheduled runner).

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$updates = @{
    "H5" = 97.666664
    "I5" = 96.35714
    "J5" = 102.25
    "K5" = 96.35714
    "L5" = 102.25
    "M5" = 18.64286
    "N5" = -332.25
    "H15" = 177.02
    "I15" = 177.02
    "K15" = 531.0600000000001
    "M15" = -362.0600000000001
    "H40" = 1828.9375
    "I40" = 1512.35
    "K40" = 1512.35
    "M40" = -1337.35
    "H116" = 3989.2942
    "I116" = 4097.3335
    "J116" = 3730
    "K116" = 4097.3335
    "L116" = 3730
    "M116" = -655.3334999999997
    "N116" = -10614
    "H137" = 939.48
    "I137" = 904.2
    "J137" = 1080.6
    "K137" = 2712.6
    "L137" = 3241.8
    "M137" = -162.6000000000004
    "N137" = -8341.799999999999
    "H138" = 2592.889
    "I138" = 985.75
    "J138" = 2999.7595
    "K138" = 2957.25
    "L138" = 8999.2785
    "M138" = 2182.75
    "N138" = -19279.2785
}
foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$updates = @{
    "H45" = 2101.1428
    "I45" = 1201.7142
    "J45" = 3000.5715
    "K45" = 1201.7142
    "L45" = 3000.5715
    "M45" = -824.7141999999999
    "N45" = -3754.5715
    "H74" = 935.2889
    "I74" = 931.9535
    "J74" = 1007
    "K74" = 931.9535
    "L74" = 1007
    "M74" = -57.95349999999996
    "N74" = -2755
    "H77" = 935.2889
    "I77" = 931.9535
    "J77" = 1007
    "K77" = 4659.7675
    "L77" = 5035
    "M77" = -291.7674999999999
    "N77" = -13771
    "H122" = 603.6
    "I122" = 611
    "J122" = 500
    "K122" = 1833
    "L122" = 1500
    "M122" = 617
    "N122" = -6400
}
foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$updates = @{
    "H9" = 20000
    "J9" = 20000
    "L9" = 20000
    "N9" = -20336
    "H22" = 369.54544
    "I22" = 354.21054
    "J22" = 466.66666
    "K22" = 354.21054
    "L22" = 466.66666
    "M22" = -181.21054
    "N22" = -812.66666
    "H86" = 2300.7827
    "I86" = 2288
    "J86" = 2312.5
    "K86" = 2288
    "L86" = 2312.5
    "M86" = -1165
    "N86" = -4558.5
    "H89" = 2300.7827
    "I89" = 2288
    "J89" = 2312.5
    "K89" = 11440
    "L89" = 11562.5
    "M89" = -5824
    "N89" = -22794.5
    "H134" = 23280.852
    "I134" = 1689.2122
    "K134" = 5067.6366
    "M134" = -2532.6366
}
foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$updates = @{
    "H58" = 3829.9412
    "I58" = 666.1053000000001
    "J58" = 7837.467
    "K58" = 666.1053000000001
    "L58" = 7837.467
    "M58" = -463.1053000000001
    "N58" = -8243.467000000001
    "H99" = 2559.5588
    "I99" = 2335.0386
    "J99" = 3289.25
    "K99" = 2335.0386
    "L99" = 3289.25
    "M99" = -837.0385999999999
    "N99" = -6285.25
    "H126" = 2559.5588
    "I126" = 2335.0386
    "J126" = 3289.25
    "K126" = 7005.1158
    "L126" = 9867.75
    "M126" = -4535.1158
    "N126" = -14807.75
    "H134" = 927.5172
    "I134" = 711.3889
    "J134" = 1281.1818
    "K134" = 2134.1667
    "L134" = 3843.5454
    "M134" = 400.8332999999998
    "N134" = -8913.545399999999
    "H136" = 3829.9412
    "I136" = 666.1053000000001
    "J136" = 7837.467
    "K136" = 1998.3159
    "L136" = 23512.401
    "M136" = 551.6840999999999
    "N136" = -28612.401
}
foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$updates = @{
    "H45" = 860.75
    "I45" = 373.33334
    "J45" = 1153.2
    "K45" = 1120.00002
    "L45" = 3459.6
    "M45" = -588.0000199999999
    "N45" = -4523.6
}
foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$updates = @{
    "H70" = 4763.9375
    "J70" = 4726.9165
    "L70" = 4726.9165
    "N70" = -5266.9165
    "H73" = 4763.9375
    "J73" = 4726.9165
    "L73" = 4726.9165
    "N73" = -6598.9165
    "H107" = 157.125
    "I107" = 157.125
    "J107" = 0
    "K107" = 157.125
    "L107" = 0
    "M107" = 1762.875
    "H122" = 11803287
    "I122" = 15965690
    "J122" = 8334616.5
    "K122" = 47897070
    "L122" = 25003849.5
    "M122" = -47894620
    "N122" = -25008749.5
    "H132" = 3074.3809
    "I132" = 3203.1667
    "J132" = 2902.6667
    "K132" = 9609.500100000001
    "L132" = 8708.000100000001
    "M132" = -7079.500100000001
    "N132" = -13768.0001
}
foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
$ws.Range("N107").ClearContents()

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$updates = @{
    "H7" = 6538383.5
    "I7" = 2709.0908
    "K7" = 2709.0908
    "M7" = -2597.0908
    "H16" = 5009.759
    "I16" = 7650
    "J16" = 2180.9285
    "K16" = 7650
    "L16" = 2180.9285
    "M16" = -7480
    "N16" = -2520.9285
    "H22" = 1756.4
    "I22" = 1931.375
    "K22" = 1931.375
    "M22" = -1636.375
    "H27" = 1756.4
    "I27" = 1931.375
    "K27" = 1931.375
    "M27" = -1824.375
    "H40" = 5051355
    "I40" = 5051355
    "J40" = 0
    "K40" = 5051355
    "L40" = 0
    "M40" = -5051219
    "H46" = 1210
    "I46" = 1227.3334
    "J46" = 1184
    "K46" = 1227.3334
    "L46" = 1184
    "M46" = -1039.3334
    "N46" = -1560
    "H122" = 21660
    "I122" = 51000
    "J122" = 2100
    "K122" = 153000
    "L122" = 6300
    "M122" = -150550
    "N122" = -11200
    "H126" = 6538383.5
    "I126" = 2709.0908
    "K126" = 8127.2724
    "M126" = -5657.2724
    "H132" = 2083.5652
    "I132" = 1904.1111
    "J132" = 2338.5789
    "K132" = 5712.3333
    "L132" = 7015.736699999999
    "M132" = -3182.3333
    "N132" = -12075.7367
}
foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
$ws.Range("N40").ClearContents()

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$updates = @{
    "H122" = 1587.3158
    "I122" = 1572.4375
    "K122" = 4717.3125
    "M122" = -2267.3125
    "H123" = 25657
    "J123" = 25657
    "L123" = 25657
    "N123" = -35457
    "H132" = 1792.6097
    "I132" = 1439.5
    "J132" = 2128.9048
    "K132" = 4318.5
    "L132" = 6386.714399999999
    "M132" = -1788.5
    "N132" = -11446.7144
    "H136" = 1222.924
    "I136" = 1148.3422
    "K136" = 3445.0266
    "M136" = -895.0266000000001
}
foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
